# faturamento_diario.xlsx edit
# - Updates the "05/2025" (May) block: corrects day 6 (row 7) and adds
#   5 new daily rows (days 7-11) right after it, pushing every subsequent
#   row (Apr/Mar/Fev blocks) down by 5 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows right after row 7 (old row 8 -> new row 13, etc.)
$ws.Rows("8:12").Insert()

# Fix the existing day-6 total for 05/2025 (row 7)
$ws.Range("B7").Value = 28248.43

# New daily rows for 05/2025 (days 7-11), filling the freshly inserted rows
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 31437.91
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 2025
$ws.Range("E8").Value = "05/2025"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 27732.15
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 2025
$ws.Range("E9").Value = "05/2025"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 25508.17
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 2025
$ws.Range("E10").Value = "05/2025"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 14802.01
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 2025
$ws.Range("E11").Value = "05/2025"

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 9716.9
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 2025
$ws.Range("E12").Value = "05/2025"
